{"js": "// Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" bullet list (\"Impact\" sub-section)\n// to use impact-focused accomplishment statements, per the commit:\n//   \"Fix Key Achievements to use proper accomplishment statements\"\n//\n// Old bullets (6) -> New bullets (4): first four bullets get new text, and the\n// final two ETL/data-warehouse bullets are removed entirely.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map of the exact original bullet text -> its replacement (per the diff).\n// The last two original bullets have no replacement (they are deleted).\nconst replacements = [\n  {\n    from: \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n    to: \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\",\n  },\n  {\n    from: \"\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n    to: \"\u2022 $4.7M savings enabled nonprofit access\",\n  },\n  {\n    from: \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n    to: \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\",\n  },\n  {\n    from: \"\u2022 Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality\",\n    to: \"\u2022 178% accuracy improvement in racial classification algorithms\",\n  },\n];\n\nconst deletions = [\n  \"\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n  \"\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\",\n];\n\n// Only touch bullets that live under the \"KEY ACHIEVEMENTS AND IMPACT\" section\n// (after that heading, before the next Heading-2 section) so we never collide\n// with visually-similar text elsewhere in the resume.\nconst items = paragraphs.items;\nlet inSection = false;\nconst toDelete = [];\n\nfor (let i = 0; i < items.length; i++) {\n  const para = items[i];\n  const text = para.text;\n\n  if (text === \"KEY ACHIEVEMENTS AND IMPACT\") {\n    inSection = true;\n    continue;\n  }\n  if (!inSection) continue;\n\n  // Stop once we reach the next top-level (Heading 2) section.\n  if (text === \"TECHNICAL SKILLS\") {\n    break;\n  }\n\n  const rep = replacements.find((r) => r.from === text);\n  if (rep) {\n    para.insertText(rep.to, Word.InsertLocation.replace);\n    continue;\n  }\n\n  if (deletions.includes(text)) {\n    toDelete.push(para);\n  }\n}\n\n// Delete the now-dropped bullets (collected first so indices/text matching\n// above aren't disturbed by the deletions).\nfor (const para of toDelete) {\n  para.delete();\n}\n\nawait context.sync();\n", "ps1": "# Rewrite the \"KEY ACHIEVEMENTS AND IMPACT\" bullet list (\"Impact\" sub-section)\n# to use impact-focused accomplishment statements, per the commit:\n#   \"Fix Key Achievements to use proper accomplishment statements\"\n#\n# Old bullets (6) -> New bullets (4): first four bullets get new text, and the\n# final two ETL/data-warehouse bullets are removed entirely.\n\n$d = $word.ActiveDocument\n\n$replacements = @{\n    \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\" = \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\";\n    \"\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\" = \"\u2022 `$4.7M savings enabled nonprofit access\";\n    \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\" = \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\";\n    \"\u2022 Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality\" = \"\u2022 178% accuracy improvement in racial classification algorithms\";\n}\n\n$deletions = @(\n    \"\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n    \"\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\"\n)\n\n# Only touch bullets that live under the \"KEY ACHIEVEMENTS AND IMPACT\" section\n# (after that heading, before the next top-level section) so we never collide\n# with visually-similar text elsewhere in the resume (e.g. the nearly\n# identical \"Trigonometric algorithm...\" bullet under Professional Experience).\n$inSection = $false\n$toDelete = New-Object System.Collections.ArrayList\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $full = $p.Range.Text\n    # Strip the trailing paragraph mark for comparison.\n    $text = $full.TrimEnd([char]13, [char]7)\n\n    if ($text -eq \"KEY ACHIEVEMENTS AND IMPACT\") {\n        $inSection = $true\n        continue\n    }\n    if (-not $inSection) {\n        continue\n    }\n    if ($text -eq \"TECHNICAL SKILLS\") {\n        break\n    }\n\n    if ($replacements.ContainsKey($text)) {\n        $p.Range.Text = $replacements[$text]\n        continue\n    }\n\n    if ($deletions -contains $text) {\n        [void]$toDelete.Add($i)\n    }\n}\n\n# Delete the now-dropped bullets by paragraph index, highest index first so\n# earlier indices collected above stay valid as we remove paragraphs.\n$sorted = $toDelete | Sort-Object -Descending\nforeach ($idx in $sorted) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
